# Add two new columns, I ("I0") and J ("IF"), to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting from an existing header cell (H1) onto the
# new header cells so they pick up the same style (bold, centered, bordered)
# without introducing a brand new style entry.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Header labels for the new columns.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New column data (I0 / IF) for each existing data row.
$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 9

$ws.Range("I3").Value = 4
$ws.Range("J3").Value = 6

$ws.Range("I4").Value = 9
$ws.Range("J4").Value = 9

$ws.Range("I5").Value = 4
$ws.Range("J5").Value = 7

$ws.Range("I6").Value = 5
$ws.Range("J6").Value = 6

$ws.Range("I7").Value = 8
$ws.Range("J7").Value = 9

$ws.Range("I8").Value = 6
$ws.Range("J8").Value = 8

$ws.Range("I9").Value = 7
$ws.Range("J9").Value = 8

$ws.Range("I10").Value = 10
$ws.Range("J10").Value = 10

$ws.Range("I11").Value = 9
$ws.Range("J11").Value = 9

$ws.Range("I12").Value = 5
$ws.Range("J12").Value = 6

$ws.Range("I13").Value = 6
$ws.Range("J13").Value = 8

$ws.Range("I14").Value = 6
$ws.Range("J14").Value = 8

$ws.Range("I15").Value = 8
$ws.Range("J15").Value = 8

$ws.Range("I16").Value = 6
$ws.Range("J16").Value = 6
